# Duplicate the existing "Mat sau" record (row 8) into three new rows (9-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Mặt sau",
    "PHAM DUY LONG",
    "S Trà Co, Thanh Cái, Qung NInh phó Móng Khu Trang Ginl Trà Co, Thanh Móng Cál, phó",
    "03/12/2006",
    "022206004066",
    "0v12/2031",
    "Việt Nam",
    "Hải Xuan, Thành phố Móng Cái, Quảng Ninh Hải Xuán, Thành phó Móng Cá",
    "Nam"
)

# Columns D (birth_day) and E (id) hold text that LOOKS like a date / a
# number with a leading zero. A plain .Value assignment would let Excel's
# auto-detection coerce them (date serial / numeric, dropping the leading
# zero), so those two columns are forced to Text format first and the
# style is reset back to Normal afterwards so the cell keeps the default
# (unstyled) appearance of the rest of the row.
$textColumns = @(4, 5)

foreach ($r in 9..11) {
    for ($c = 1; $c -le $values.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textColumns -contains $c) {
            $cell.NumberFormat = "@"
            $cell.Value = $values[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$c - 1]
        }
    }
}
